# feat: add 2022-Q1 data
$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new "2022-Q1" sheet right after "2021-Q4" (i.e. right
#    before the "总计" summary sheet) and copy the formatting used by
#    the other quarterly detail sheets.
# ------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $ws4)
$newSheet.Name = "2022-Q1"

# Copy the header-row formatting (bold, centered, bordered) and the
# row-index column formatting from the "2021-Q4" sheet so the new
# sheet matches the established look of the other quarter sheets.
$ws4.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$ws4.Range("A2").Copy()
$newSheet.Range("A2:A5").PasteSpecial(-4122)

# Header row
$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"

# Row index column (A2:A5) -- 0-based row counter, stored as numbers.
$newSheet.Cells.Item(2, 1).Value = 0
$newSheet.Cells.Item(3, 1).Value = 1
$newSheet.Cells.Item(4, 1).Value = 2
$newSheet.Cells.Item(5, 1).Value = 3

function Set-TextCell($sheet, $row, $col, $text) {
    $c = $sheet.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $text
}

# Row 2 -- 506001
Set-TextCell $newSheet 2 2 "506001"
Set-TextCell $newSheet 2 3 "万家科创板 2 年定期开放混合型证券投资基金"
Set-TextCell $newSheet 2 4 "12.84"
Set-TextCell $newSheet 2 5 "98.14"
Set-TextCell $newSheet 2 6 "3.38"
Set-TextCell $newSheet 2 7 "0.4340"
$newSheet.Cells.Item(2, 8).Value = 8

# Row 3 -- 180028
Set-TextCell $newSheet 3 2 "180028"
Set-TextCell $newSheet 3 3 "银华永祥灵活配置混合"
Set-TextCell $newSheet 3 4 "0.61"
Set-TextCell $newSheet 3 5 "77.23"
Set-TextCell $newSheet 3 6 "3.56"
Set-TextCell $newSheet 3 7 "0.0217"
$newSheet.Cells.Item(3, 8).Value = 5

# Row 4 -- 006689
Set-TextCell $newSheet 4 2 "006689"
Set-TextCell $newSheet 4 3 "方正富邦信泓灵活配置混合A"
Set-TextCell $newSheet 4 4 "0.06"
Set-TextCell $newSheet 4 5 "89.81"
Set-TextCell $newSheet 4 6 "4.76"
Set-TextCell $newSheet 4 7 "0.0029"
$newSheet.Cells.Item(4, 8).Value = 7

# Row 5 -- 008182
Set-TextCell $newSheet 5 2 "008182"
Set-TextCell $newSheet 5 3 "方正富邦信泓灵活配置混合C"
Set-TextCell $newSheet 5 4 "0.04"
Set-TextCell $newSheet 5 5 "89.81"
Set-TextCell $newSheet 5 6 "4.76"
Set-TextCell $newSheet 5 7 "0.0019"
$newSheet.Cells.Item(5, 8).Value = 7

# ------------------------------------------------------------------
# 2. Update the "总计" (Total) summary sheet: insert a new data row
#    at the top (row 2) for "2022-Q1" and push the rest down.
# ------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Rows.Item(2).Insert()

# Inserting a row copies the formatting of the row above (the bold
# header row) onto the new row -- clear that back to the plain,
# unstyled look the data rows use, then restore just the row-index
# column's style (A3 still carries it) to match A3:A7 below.
$wsTotal.Range("A2:D2").ClearFormats()
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q1"
$wsTotal.Cells.Item(2, 3).Value = 4
$wsTotal.Cells.Item(2, 4).Value = 0.46

# Re-number the 0-based row-index column (A) for the rows that moved
# down, so it continues 0,1,2,3,4,5 as before.
$wsTotal.Cells.Item(3, 1).Value = 1
$wsTotal.Cells.Item(4, 1).Value = 2
$wsTotal.Cells.Item(5, 1).Value = 3
$wsTotal.Cells.Item(6, 1).Value = 4
$wsTotal.Cells.Item(7, 1).Value = 5
